# Actualización automática 2025-08-07 16:40:08
# Apply updated sales figure (806.01) for client "FUENTES PAREDES MARIA FERNANDA"
# (group PIEDRA SINTERIZADA, month agosto) and propagate it through the
# related summary/aggregate sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO -------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("L12").Value = 806.01
$wsGrupo.Range("L33").Value = "2 de 31"

# --- Sheet: VENTA MENSUAL ----------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F12").Value = 1263.93
$wsMensual.Range("F33").Value = 4530.17

# --- Sheet: CUMPLIMIENTO MENSUAL ---------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Column E width shrinks slightly (stored width 23 -> 22)
$wsCumpl.Columns.Item(5).ColumnWidth = 21.17

# Row 15 (PIEDRA SINTERIZADA)
$wsCumpl.Range("D15").Value = 1027.73
$wsCumpl.Range("E15").Value = -500.7
$wsCumpl.Range("F15").Value = 1.950040794641671

# Row 19 (TOTAL)
$wsCumpl.Range("D19").Value = 4530.169999999999
$wsCumpl.Range("E19").Value = 27579.11107555787
$wsCumpl.Range("F19").Value = 0.1410859990711048
